$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "campaign_fFLML21TlvlGAzHlw3SLlHv0difPft"
$ws.Range("B2").Value = "INACTIVE"
$ws.Range("C2").Value = "play"
$ws.Range("D2").Value = 4

# Row 3
$ws.Range("A3").Value = "campaign_jCiEsmQGldMt4a2nER9k3N5zkpVwb"
$ws.Range("B3").Value = "INACTIVE"
$ws.Range("C3").Value = "mouse"
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 0
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "0.0"
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0

# Row 4
$ws.Range("A4").Value = "campaign_Hk5mBwb0NxvbzVPiCvlDw7ATXaCAs3"
$ws.Range("C4").Value = "roladin"
$ws.Range("D4").Value = 6
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "0.0"
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
